$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.9184, 0.9307, 0.9076),
    @(0.9152, 0.9256, 0.911),
    @(0.0058, 0.0087, 0.0073),
    @(0.9488, 0.953, 0.9392),
    @(0.9471, 0.9504, 0.9408),
    @(0.0033, 0.0051, 0.0039),
    @(10.9161, 7.842, 25.4019),
    @(5.0981, 0.8906, 0.8068)
)

$row = 2
foreach ($values in $data) {
    $ws.Range("D$row").Value = $values[0]
    $ws.Range("E$row").Value = $values[1]
    $ws.Range("F$row").Value = $values[2]
    $row++
}
